$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking strings
# (e.g. "210.23") are preserved as text rather than being converted to numbers.
$ws.Range('D2,D3,D5,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D20,D22,D24,D25,D27,D32,D33,D34,D39,D40,D43,D44,D45,D46,D47,D49,D50,D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.304.44'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '1.591.55'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = '210.23'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').Value = '0.0611'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').Value = '0.245'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '1.814.22'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.586.91'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '4.07'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '0.518'
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '26.328.76'
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('D20').Value = '212.08'
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D24').Value = '8.92'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').Value = '145.25'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = '7.06'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').Value = '2.98'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').Value = '1.300.49'
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E36').Value = '  -2.09%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '1.10'
$ws.Range('E39').Value = '  -13.59%  '
$ws.Range('D40').Value = '0.814'
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.14'
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '62.70'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '0.762'
$ws.Range('E45').Value = '  -1.98%  '
$ws.Range('D46').Value = '1.727.42'
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('D47').Value = '88.48'
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('E48').Value = '  -3.86%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -3.29%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0988'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0506'
$ws.Range('E51').Value = '  -1.34%  '
